$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 100
$ws.Range("D5").Value = 100
